# Simulated Wild Card round and logged it
$wb = $excel.ActiveWorkbook

# Update "OFF" sheet (sheet1.xml) - Row 3 ("R")
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 391
$wsOff.Range("C3").Value = 263
$wsOff.Range("D3").Value = 93
$wsOff.Range("E3").Value = 31

# Update "DEF" sheet (sheet2.xml) - Row 3 ("R")
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 514
$wsDef.Range("C3").Value = 365
$wsDef.Range("D3").Value = 110
$wsDef.Range("E3").Value = 56
